$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update noise dB SPL values (accounting for frequency response curve of Sennheiser HD280)
$ws.Range("B2").Value = 80
$ws.Range("C2").Formula = "=0.4*B2"
$ws.Range("I2").Value = 100
$ws.Range("K2").Value = 100

# Update the active selection on the sheet
$ws.Range("K12").Select()
